# Updated symbol list on Thu Jan 19 20:53:09 UTC 2023 with GitHub Actions
#
# Refresh the "Price" (D) and "Volume(1h)" (E) columns on the crypto
# tracker sheet with the latest scraped quotes. Values are plain text
# (prices keep their original decimal precision, volume keeps its "%"
# suffix) so each cell is forced to text storage via NumberFormat "@"
# before the write, then restored to the sheet's default "Normal"
# style so no stray formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 cell -> new text value
$updates = [ordered]@{
    "D2"  = "294.57";    "E2"  = "1.31%"
    "D3"  = "31.01";     "E3"  = "0.41%"
    "D4"  = "4.964";     "E4"  = "2.23%"
    "D5"  = "0.07339";   "E5"  = "2.42%"
    "D6"  = "2.311";     "E6"  = "33.52%"
    "D7"  = "7.726";     "E7"  = "0.77%"
    "D8"  = "3.744";     "E8"  = "-0.36%"
    "D9"  = "0.9083";    "E9"  = "1.36%"
                          "E10" = "2.13%"
    "D11" = "0.08013";   "E11" = "8.48%"
    "D12" = "0.08058";   "E12" = "0.43%"
    "D13" = "0.03104";   "E13" = "2.33%"
    "D14" = "0.1008";    "E14" = "1.06%"
    "D15" = "0.001520";  "E15" = "1.60%"
    "D16" = "0.005803";  "E16" = "2.46%"
    "D17" = "3.490";     "E17" = "0.84%"
    "D18" = "2.076";     "E18" = "-1.62%"
    "D19" = "0.3327";    "E19" = "1.08%"
    "D20" = "0.1304";    "E20" = "0.11%"
    "D21" = "3.962";     "E21" = "-9.72%"
    "D22" = "0.2101";    "E22" = "4.70%"
    "D23" = "0.04545";   "E23" = "1.12%"
    "D24" = "0.001213";  "E24" = "-0.36%"
    "D25" = "0.004652";  "E25" = "16.09%"
                          "E26" = "3.37%"
    "D27" = "0.0003396"; "E27" = "-95.47%"
    "D39" = "0.01606";   "E39" = "-1.82%"
    "D40" = "0.04427";   "E40" = "2.20%"
    "D41" = "0.007398";  "E41" = "-0.14%"
                          "E42" = "1.95%"
    "D43" = "0.008609"
    "D44" = "0.002048";  "E44" = "0.02%"
    "D45" = "0.009518";  "E45" = "-14.27%"
    "D46" = "0.00005924";"E46" = "3.22%"
                          "E47" = "-0.43%"
                          "E48" = "2.69%"
    "D49" = "0.002900";  "E49" = "-3.70%"
    "D50" = "0.00002101";"E50" = "-0.43%"
    "D51" = "0.0002001"; "E51" = "-0.43%"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text storage so the numeric-looking / percent strings keep
    # their exact original representation instead of being parsed into
    # a Double (and reformatted) by the normal Value auto-detection.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    # The source cells carry no explicit style (default "Normal"); put
    # it back so only the text content changes, not the formatting.
    $cell.Style = "Normal"
}
